$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values per repull of data / mean calculation fix
$ws.Range("F2").Value = 3
$ws.Range("F3").Value = -4
$ws.Range("F4").Value = -8
$ws.Range("F5").Value = -7
$ws.Range("F6").Value = 3
$ws.Range("F12").Value = -1
